$d = $word.ActiveDocument

$replacements = @(
    @{old = "353×2=706";   new = "675×5=3375"},
    @{old = "955×2=1910";  new = "414×5=2070"},
    @{old = "430×4=1720";  new = "668×7=4676"},
    @{old = "634×8=5072";  new = "744×8=5952"},
    @{old = "975×4=3900";  new = "799×6=4794"},
    @{old = "987×4=3948";  new = "210×2=420"},
    @{old = "324×4=1296";  new = "350×9=3150"},
    @{old = "296×7=2072";  new = "936×6=5616"},
    @{old = "592×2=1184";  new = "156×6=936"},
    @{old = "265×2=530";   new = "659×3=1977"},
    @{old = "255×8=2040";  new = "493×4=1972"},
    @{old = "446×4=1784";  new = "706×6=4236"},
    @{old = "916×9=8244";  new = "231×8=1848"},
    @{old = "268×9=2412";  new = "236×4=944"},
    @{old = "488×3=1464";  new = "279×3=837"},
    @{old = "524×4=2096";  new = "339×2=678"},
    @{old = "907×3=2721";  new = "716×4=2864"},
    @{old = "688×5=3440";  new = "877×3=2631"},
    @{old = "336×2=672";   new = "630×7=4410"},
    @{old = "243×8=1944";  new = "807×3=2421"},
    @{old = "492×4=1968";  new = "747×9=6723"},
    @{old = "648×5=3240";  new = "333×2=666"},
    @{old = "410×7=2870";  new = "760×2=1520"},
    @{old = "699×8=5592";  new = "377×9=3393"},
    @{old = "981×4=3924";  new = "593×5=2965"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
